$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values (Date serial, US_Gasoline, UK_Gasoline, Nigeria_Gasoline),
# sorted by Date descending, written in place over A2:D77 (row-by-row so
# other cells/attributes on each row, e.g. H3, are left untouched).
$data = @(
    @(44089.0,2.274,5.742620935843709,1.5771921323985885),
    @(44058.0,2.272,5.719967515887019,1.4487530851112151),
    @(44027.0,2.272,5.637896355685364,1.419226867344003),
    @(43997.0,2.17,5.368416783434191,1.2460063897763576),
    @(43966.0,1.961,5.314893709825982,1.2154959647502384),
    @(43936.0,1.938,5.527464439067503,1.2302590736338446),
    @(43905.0,2.329,6.099077884179776,1.4271005254152596),
    @(43876.0,2.533,6.268388641662425,1.413977761963165),
    @(43845.0,2.636,6.4491436222269165,1.419226867344003),
    @(43814.0,2.645,6.310883322779865,1.413977761963165),
    @(43784.0,2.693,6.3732990938736975,1.419226867344003),
    @(43753.0,2.724,6.445495990085974,1.4123374165316536),
    @(43723.0,2.681,6.441738250254035,1.4090567256686333),
    @(43692.0,2.707,6.518591461676473,1.430348409369653),
    @(43661.0,2.823,6.461515651201268,1.407416380237118),
    @(43631.0,2.804,6.4739844617209314,1.419226867344003),
    @(43600.0,2.946,6.496258240429763,1.4172584528261887),
    @(43570.0,2.881,6.294687993553611,1.4238198345522357),
    @(43539.0,2.594,6.107836473870946,1.4172584528261887),
    @(43511.0,2.393,6.028862567750117,1.4221794891207242),
    @(43480.0,2.338,6.0593768068763385,1.419226867344003),
    @(43449.0,2.457,6.136302713966084,1.426608421785806),
    @(43419.0,2.736,6.523736966977912,1.4211952818618172),
    @(43388.0,2.943,6.638905368177693,1.4227700134760684),
    @(43358.0,2.915,6.632295151831273,1.4172584528261887),
    @(43327.0,2.914,6.523989720953724,1.4231636963796312),
    @(43296.0,2.928,6.473354436820707,1.4271005254152596),
    @(43266.0,2.97,6.489948661458384,1.4221794891207242),
    @(43235.0,2.987,6.323802855304136,1.4267068425116969),
    @(43205.0,2.873,6.116060411440623,1.426608421785806),
    @(43174.0,2.709,6.041765080404269,1.4271005254152596),
    @(43146.0,2.705,6.1600749805786394,1.4271005254152596),
    @(43115.0,2.671,6.145842149646608,1.574731614251321),
    @(43084.0,2.594,6.0866369019199125,1.4421917033851681),
    @(43054.0,2.678,6.042552324750687,1.4172584528261887),
    @(43023.0,2.621,5.942380613089769,1.4172584528261887),
    @(42993.0,2.761,6.032861502767342,1.3975743076480474),
    @(42962.0,2.494,5.865818058807853,1.4123374165316536),
    @(42931.0,2.414,5.777753969010948,1.4271005254152596),
    @(42901.0,2.46,5.861139236005651,1.4271005254152596),
    @(42870.0,2.503,5.859758157349301,1.4271005254152596),
    @(42840.0,2.528,5.950069267601271,1.4271005254152596),
    @(42809.0,2.437,6.056001393031812,1.4271005254152596),
    @(42781.0,2.416,6.0799684106290535,1.4271005254152596),
    @(42750.0,2.458,6.020746932066157,1.419226867344003),
    @(42719.0,2.366,5.78626760948217,1.4271005254152596),
    @(42689.0,2.295,5.878182312696757,1.4231636963796312),
    @(42658.0,2.359,5.760030059306622,1.4398952197810517),
    @(42628.0,2.327,5.641130535060944,1.4103690020138395),
    @(42597.0,2.284,5.531489688130462,1.4271005254152596),
    @(42566.0,2.345,5.664048165106963,1.4271005254152596),
    @(42536.0,2.467,5.628566796507543,1.4271005254152596),
    @(42505.0,2.371,5.5002691658465395,1.4271005254152596),
    @(42475.0,2.216,5.399263098093157,1.053101767030571),
    @(42444.0,2.071,5.160046895616645,1.0324334145935223),
    @(42415.0,1.872,5.14359585218118,0.912360129006859),
    @(42384.0,2.057,5.160834526930455,0.9802704298714472),
    @(42353.0,2.144,5.259088468317835,1.1296730917735414),
    @(42323.0,2.26,5.4396470624737505,1.0112729585270201),
    @(42292.0,2.387,5.523706230853578,1.0649122541374558),
    @(42262.0,2.462,5.655438187935245,0.9113759217479519),
    @(42231.0,2.726,5.807064765751899,0.9310600669260934),
    @(42200.0,2.88,5.904502377874318,0.925154823372651),
    @(42170.0,2.885,5.904147502512213,0.957141559287131),
    @(42139.0,2.802,5.8713416863346595,0.9620625955816664),
    @(42109.0,2.555,5.708917452151156,0.9212179943370227),
    @(42078.0,2.546,5.632555712522304,0.9620625955816664),
    @(42050.0,2.301,5.437428741017294,0.8808654967218327),
    @(42019.0,2.208,5.500826016751238,0.9546810411398633),
    @(41988.0,2.632,5.895308019721308,0.9546810411398633),
    @(41958.0,2.997,6.212591486152163,0.9546810411398633),
    @(41927.0,3.255,6.429726863232947,0.9546810411398633),
    @(41897.0,3.484,6.518793700672069,0.9546810411398633),
    @(41866.0,3.565,6.55715902896345,0.9546810411398633),
    @(41835.0,3.688,6.651015999128547,0.9546810411398633),
    @(41805.0,3.766,6.57890989961086,0.9546810411398633)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Auto-fit columns C:D to their content
$ws.Columns.Item(3).ColumnWidth = 11.1
$ws.Columns.Item(4).ColumnWidth = 14

# Restore the last-used selection recorded after the edit
$ws.Range("G17").Select()
